# Update "想去人数" (want-to-go headcount, column F) figures on the
# "展览" and "全部类型" sheets to the refreshed counts from the latest
# scrape (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# row => new value, per sheet (only column F changes)
$updates1 = @{
    4  = 104
    5  = 1714
    7  = 933
    9  = 2044
    10 = 1054
    11 = 566
    17 = 79
    18 = 125
    19 = 1498
    20 = 563
    21 = 665
    22 = 549
    23 = 11940
    24 = 11950
    27 = 267
    28 = 1874
}

$updates4 = @{
    6  = 104
    7  = 1714
    9  = 933
    11 = 2044
    12 = 1054
    13 = 566
    20 = 79
    22 = 125
    23 = 1498
    24 = 563
    25 = 665
    26 = 549
    27 = 11940
    28 = 11950
    31 = 267
    32 = 1874
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates1.Keys) {
    $ws1.Range("F$row").Value = $updates1[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}
